# Updates cryptos list prices/volumes to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.037.08'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '1.830.74'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'241.47"
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = "'0.6538"
$ws.Range('E6').Value = '  -3.11%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = "'44.59"
$ws.Range('E8').Value = '  +5.75%  '
$ws.Range('D9').Value = "'0.07364"
$ws.Range('E9').Value = '  -1.25%  '
$ws.Range('D10').Value = "'0.2940"
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('D11').Value = "'22.98"
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').Value = "'0.07674"
$ws.Range('E12').Value = '  -0.59%  '
$ws.Range('D13').Value = '1.833.39'
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = "'82.35"
$ws.Range('E16').Value = '  -4.61%  '
$ws.Range('D17').Value = "'6.074"
$ws.Range('E17').Value = '  -1.55%  '
$ws.Range('E18').Value = '  +3.31%  '
$ws.Range('D19').Value = '29.031.60'
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('D20').Value = '2.082.90'
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('E21').Value = '  -0.96%  '
$ws.Range('D22').Value = "'224.36"
$ws.Range('E22').Value = '  -2.11%  '
$ws.Range('D23').Value = "'0.9998"
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').Value = "'7.114"
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = "'158.24"
$ws.Range('E26').Value = '  -1.74%  '
$ws.Range('D27').Value = "'8.520"
$ws.Range('D28').Value = "'0.1383"
$ws.Range('E28').Value = '  -1.95%  '
$ws.Range('D29').Value = "'17.93"
$ws.Range('E29').Value = '  -0.60%  '
$ws.Range('D30').Value = "'1.502"
$ws.Range('E30').Value = '  -0.72%  '
$ws.Range('D31').Value = "'4.114"
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('D33').Value = "'4.014"
$ws.Range('E33').Value = '  -1.57%  '
$ws.Range('D34').Value = "'0.05336"
$ws.Range('E34').Value = '  +0.40%  '
$ws.Range('D35').Value = "'0.7428"
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('D36').Value = "'1.834"
$ws.Range('E36').Value = '  -2.38%  '
$ws.Range('D37').Value = "'1.153"
$ws.Range('E37').Value = '  +1.17%  '
$ws.Range('E38').Value = '  -1.17%  '
$ws.Range('D39').Value = '1.291.68'
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('D40').Value = "'0.01787"
$ws.Range('E40').Value = '  -1.23%  '
$ws.Range('D41').Value = "'2.747"
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('D42').Value = "'6.361"
$ws.Range('E42').Value = '  +6.34%  '
$ws.Range('D43').Value = "'0.8941"
$ws.Range('E43').Value = '  -2.75%  '
$ws.Range('D44').Value = "'0.9990"
$ws.Range('D45').Value = "'103.11"
$ws.Range('E45').Value = '  -0.52%  '
$ws.Range('D46').Value = '1.980.50'
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('D47').Value = "'0.5140"
$ws.Range('D48').Value = "'64.20"
$ws.Range('E48').Value = '  -0.20%  '
$ws.Range('E49').Value = '  -0.80%  '
$ws.Range('D50').Value = "'1.733"
$ws.Range('E50').Value = '  -2.85%  '
$ws.Range('D51').Value = "'0.07517"
$ws.Range('E51').Value = '  -8.02%  '
